$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2033.1666
$ws.Range("I2").Value = 2079.8
$ws.Range("K2").Value = 2079.8
$ws.Range("M2").Value = -1966.8
$ws.Range("H28").Value = 882.2857
$ws.Range("I28").Value = 578.0769
$ws.Range("K28").Value = 578.0769
$ws.Range("M28").Value = -93.07690000000002
$ws.Range("H32").Value = 2625.5217
$ws.Range("J32").Value = 2709.524
$ws.Range("L32").Value = 2709.524
$ws.Range("N32").Value = -3361.524
$ws.Range("H43").Value = 135569.38
$ws.Range("I43").Value = 7990
$ws.Range("J43").Value = 153795
$ws.Range("K43").Value = 7990
$ws.Range("L43").Value = 153795
$ws.Range("M43").Value = -7921
$ws.Range("N43").Value = -153933
$ws.Range("H134").Value = 77498.336
$ws.Range("I134").Value = 40000
$ws.Range("J134").Value = 84998
$ws.Range("K134").Value = 40000
$ws.Range("L134").Value = 84998
$ws.Range("M134").Value = -34930
$ws.Range("N134").Value = -95138

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 31750
$ws.Range("I37").Value = 23500
$ws.Range("J37").Value = 40000
$ws.Range("K37").Value = 23500
$ws.Range("L37").Value = 40000
$ws.Range("M37").Value = -23227
$ws.Range("N37").Value = -40546
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H102").Value = 6343.7617
$ws.Range("I102").Value = 5906.316
$ws.Range("K102").Value = 5906.316
$ws.Range("M102").Value = -4284.316
$ws.Range("H132").Value = 53389.87
$ws.Range("I132").Value = 2427.439
$ws.Range("J132").Value = 401633.16
$ws.Range("K132").Value = 7282.316999999999
$ws.Range("L132").Value = 1204899.48
$ws.Range("M132").Value = -4752.316999999999
$ws.Range("N132").Value = -1209959.48

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 92705
$ws.Range("I97").Value = 21569.334
$ws.Range("J97").Value = 519519
$ws.Range("K97").Value = 21569.334
$ws.Range("L97").Value = 519519
$ws.Range("M97").Value = -20578.334
$ws.Range("N97").Value = -521501
$ws.Range("H99").Value = 3262.7693
$ws.Range("I99").Value = 2563.6
$ws.Range("K99").Value = 2563.6
$ws.Range("M99").Value = -1065.6
$ws.Range("H132").Value = 146166.17
$ws.Range("J132").Value = 150399.4
$ws.Range("L132").Value = 150399.4
$ws.Range("N132").Value = -160519.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5374.4287
$ws.Range("I31").Value = 3414.6843
$ws.Range("K31").Value = 3414.6843
$ws.Range("M31").Value = -3119.6843
$ws.Range("H34").Value = 5374.4287
$ws.Range("I34").Value = 3414.6843
$ws.Range("K34").Value = 3414.6843
$ws.Range("M34").Value = -3212.6843
$ws.Range("H52").Value = 86962
$ws.Range("J52").Value = 86962
$ws.Range("L52").Value = 86962
$ws.Range("N52").Value = -87550
$ws.Range("H99").Value = 5445.1665
$ws.Range("J99").Value = 3249.5
$ws.Range("L99").Value = 3249.5
$ws.Range("N99").Value = -6245.5
$ws.Range("H104").Value = 54969.832
$ws.Range("J104").Value = 54969.832
$ws.Range("L104").Value = 54969.832
$ws.Range("N104").Value = -60211.832
$ws.Range("H126").Value = 5445.1665
$ws.Range("J126").Value = 3249.5
$ws.Range("L126").Value = 9748.5
$ws.Range("N126").Value = -14688.5
$ws.Range("H137").Value = 142362.5
$ws.Range("J137").Value = 149816.67
$ws.Range("L137").Value = 149816.67
$ws.Range("N137").Value = -160016.67
$ws.Range("H139").Value = 84993
$ws.Range("J139").Value = 84993
$ws.Range("L139").Value = 84993
$ws.Range("N139").Value = -95273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 931.1667
$ws.Range("I107").Value = 331.75
$ws.Range("J107").Value = 1230.875
$ws.Range("K107").Value = 995.25
$ws.Range("L107").Value = 3692.625
$ws.Range("M107").Value = 924.75
$ws.Range("N107").Value = -7532.625
$ws.Range("H127").Value = 3766.6667
$ws.Range("J127").Value = 3766.6667
$ws.Range("L127").Value = 11300.0001
$ws.Range("N127").Value = -21220.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 50420
$ws.Range("J106").Value = 50420
$ws.Range("L106").Value = 50420
$ws.Range("N106").Value = -52944
$ws.Range("H126").Value = 2245.45
$ws.Range("I126").Value = 1950.6875
$ws.Range("J126").Value = 3424.5
$ws.Range("K126").Value = 5852.0625
$ws.Range("L126").Value = 10273.5
$ws.Range("M126").Value = -3382.0625
$ws.Range("N126").Value = -15213.5
$ws.Range("H131").Value = 63499.5
$ws.Range("J131").Value = 63499.5
$ws.Range("L131").Value = 63499.5
$ws.Range("N131").Value = -73579.5
$ws.Range("H132").Value = 2650
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
$ws.Range("H134").Value = 85856
$ws.Range("J134").Value = 85856
$ws.Range("L134").Value = 257568
$ws.Range("N134").Value = -262638
$ws.Range("H136").Value = 58883.668
$ws.Range("J136").Value = 58883.668
$ws.Range("L136").Value = 176651.004
$ws.Range("N136").Value = -181751.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5556.3076
$ws.Range("I132").Value = 2382.6316
$ws.Range("J132").Value = 14170.571
$ws.Range("K132").Value = 7147.8948
$ws.Range("L132").Value = 42511.713
$ws.Range("M132").Value = -4617.8948
$ws.Range("N132").Value = -47571.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1357.3334
$ws.Range("I107").Value = 1397.7778
$ws.Range("J107").Value = 1316.8889
$ws.Range("K107").Value = 4193.3334
$ws.Range("L107").Value = 3950.6667
$ws.Range("M107").Value = -2273.3334
$ws.Range("N107").Value = -7790.6667
$ws.Range("H116").Value = 45680
$ws.Range("J116").Value = 45680
$ws.Range("L116").Value = 45680
$ws.Range("N116").Value = -54858
$ws.Range("H126").Value = 5379.8
$ws.Range("I126").Value = 4224.75
$ws.Range("K126").Value = 12674.25
$ws.Range("M126").Value = -10204.25
$ws.Range("H132").Value = 1315.6522
$ws.Range("I132").Value = 1066.4736
$ws.Range("J132").Value = 2499.25
$ws.Range("K132").Value = 3199.4208
$ws.Range("L132").Value = 7497.75
$ws.Range("M132").Value = -669.4207999999999
$ws.Range("N132").Value = -12557.75

Write-Host "Applied 162 cell updates"